$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.149.13"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "1.577.83"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.18"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.496"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -3.33%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0608"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.50"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "1.800.11"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").Value = "1.601.66"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.42"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "26.146.46"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.24"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "207.89"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.88"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.92"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.96"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.20"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0505"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").Value = "1.277.75"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.46"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.607"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.11"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -5.50%  "
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.56"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +2.90%  "
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.764"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.35"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").Value = "1.713.20"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("E48").Value = "  -1.75%  "
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("E51").Value = "  +0.06%  "
